# Auto-generated edit script: updates crypto price/volume columns (D, E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text values (not ambiguous with numbers) ---
$ws.Range("D2").Value = "66.756.69"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "3.090.66"
$ws.Range("E3").Value = "  +5.26%  "
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("E6").Value = "  +6.16%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.086.13"
$ws.Range("E8").Value = "  +5.23%  "
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("E11").Value = "  +4.19%  "
$ws.Range("E12").Value = "  +4.94%  "
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("E14").Value = "  +6.09%  "
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "3.601.89"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "66.734.24"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "3.090.59"
$ws.Range("E19").Value = "  +5.19%  "
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("E21").Value = "  +5.35%  "
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("E23").Value = "  +3.20%  "
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("E25").Value = "  +4.29%  "
$ws.Range("E26").Value = "  +8.48%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("E31").Value = "  +4.07%  "
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +3.09%  "
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("E38").Value = "  +6.03%  "
$ws.Range("E39").Value = "  +6.02%  "
$ws.Range("E40").Value = "  +6.08%  "
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "2.790.97"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +6.93%  "
$ws.Range("E51").Value = "  +1.69%  "

# --- Values that look like plain numbers: assign via a quoted-string formula,
#     then Copy + PasteSpecial(xlPasteValues) to collapse to a literal text value
#     without Excel auto-converting the text to a numeric cell. ---
$ws.Range("D5").Formula = '="580.29"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="168.60"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D14").Formula = '="36.47"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="16.27"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="468.12"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D23").Formula = '="7.51"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D24").Formula = '="84.11"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D26").Formula = '="13.12"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D27").Formula = '="10.11"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D30").Formula = '="2.40"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D33").Formula = '="28.33"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D35").Formula = '="1.00"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="47.28"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="50.32"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D43").Formula = '="8.70"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("D44").Formula = '="2.82"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D46").Formula = '="382.92"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("D48").Formula = '="135.04"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("D50").Formula = '="24.94"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)

$excel.CutCopyMode = $false
